$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.006.38"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "3.405.54"
$ws.Range("E3").Value = "  -2.40%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "407.66"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.95"
$ws.Range("E6").Value = "  +3.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.591"
$ws.Range("E7").Value = "  -1.36%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("E10").Value = "  -6.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.58"
$ws.Range("E11").Value = "  -1.82%  "
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("E13").Value = "  -4.02%  "
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").Value = "3.398.89"
$ws.Range("E15").Value = "  -3.29%  "
$ws.Range("D16").Value = "61.963.77"
$ws.Range("E16").Value = "  -1.39%  "
$ws.Range("E17").Value = "  -3.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.00"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("E19").Value = "  -5.61%  "
$ws.Range("E20").Value = "  -5.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "84.54"
$ws.Range("E21").Value = "  +2.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "313.93"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.91"
$ws.Range("E23").Value = "  -2.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.15"
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.74"
$ws.Range("E25").Value = "  +8.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.63"
$ws.Range("E26").Value = "  -2.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.24"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.83"
$ws.Range("E28").Value = "  +4.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.59"
$ws.Range("E29").Value = "  -3.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.175"
$ws.Range("E30").Value = "  -3.73%  "
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "42.47"
$ws.Range("E32").Value = "  -3.95%  "
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.35"
$ws.Range("E34").Value = "  -6.72%  "
$ws.Range("E35").Value = "  -2.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.74"
$ws.Range("E36").Value = "  -1.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("E38").Value = "  -5.09%  "
$ws.Range("E39").Value = "  -3.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.302"
$ws.Range("E40").Value = "  +4.41%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.125"
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "137.34"
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("E43").Value = "  -1.40%  "
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.71"
$ws.Range("E45").Value = "  -6.51%  "
$ws.Range("E46").Value = "  -1.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.38"
$ws.Range("E47").Value = "  -5.27%  "
$ws.Range("D48").Value = "2.118.22"
$ws.Range("E48").Value = "  -4.61%  "
$ws.Range("E49").Value = "  -2.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.90"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.65"
$ws.Range("E51").Value = "  +16.13%  "
